# New crime data collected - weekly CompStat update for the 84th Precinct
# Moves the report window forward one week (6/3-6/9/2024 -> 6/10-6/16/2024,
# issue number 23 -> 24) and refreshes every Week-to-Date / 28-Day /
# Year-to-Date / 2-Year crime statistic in the table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: volume/issue number and the reporting week dates.
# ---------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 31   Number  24"
$ws.Range("C9").Value  = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# ---------------------------------------------------------------------
# Row 14 (Murder): M14 flips from the "***.*" placeholder text to a
# real numeric 0, picking up the existing "#,##0.0" style (s=16).
# ---------------------------------------------------------------------
$ws.Range("M14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M14").Value = 0

# ---------------------------------------------------------------------
# Row 15 (Rape): C15 flips from a numeric 4 to the "0" placeholder text
# (shared with the other text-style cells in this row), reusing the
# existing text style (s=14) by copying it from an already-text cell.
# ---------------------------------------------------------------------
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("N15").Value = -50

# ---------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 80
$ws.Range("I16").Value = 119
$ws.Range("J16").Value = 61
$ws.Range("K16").Value = 95.081967213114
$ws.Range("L16").Value = 88.888888888888
$ws.Range("M16").Value = 26.595744680851
$ws.Range("N16").Value = -77.881040892193

# ---------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 69.230769230769
$ws.Range("I17").Value = 140
$ws.Range("J17").Value = 106
$ws.Range("K17").Value = 32.075471698113
$ws.Range("L17").Value = 91.780821917808
$ws.Range("M17").Value = 164.150943396226
$ws.Range("N17").Value = -25.133689839572

# ---------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 28
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 99
$ws.Range("J18").Value = 107
$ws.Range("K18").Value = -7.476635514018
$ws.Range("L18").Value = 2.061855670103
$ws.Range("M18").Value = 94.117647058823
$ws.Range("N18").Value = -71.304347826087

# ---------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 18.181818181818
$ws.Range("I19").Value = 325
$ws.Range("J19").Value = 312
$ws.Range("K19").Value = 4.166666666666
$ws.Range("L19").Value = 17.753623188405
$ws.Range("M19").Value = 64.974619289340
$ws.Range("N19").Value = -27.616926503340

# ---------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 28.571428571428
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 33
$ws.Range("K20").Value = -6.060606060606
$ws.Range("L20").Value = 19.230769230769
$ws.Range("M20").Value = 47.619047619047
$ws.Range("N20").Value = -90.342679127725

# ---------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 22.222222222222
$ws.Range("F21").Value = 155
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = 40.909090909090
$ws.Range("I21").Value = 719
$ws.Range("J21").Value = 622
$ws.Range("K21").Value = 15.594855305466
$ws.Range("L21").Value = 33.148148148148
$ws.Range("M21").Value = 72.009569377990
$ws.Range("N21").Value = -61.177105831533

# ---------------------------------------------------------------------
# Row 22 (Transit): C22 flips from the "0" placeholder text to a real
# numeric 1, picking up the existing "#,##0" style (s=15).
# ---------------------------------------------------------------------
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 18
$ws.Range("K22").Value = -14.285714285714
$ws.Range("L22").Value = -18.181818181818
$ws.Range("M22").Value = -45.454545454545

# ---------------------------------------------------------------------
# Row 23 (Housing): C23 flips from a numeric 1 back to the "0"
# placeholder text, reusing the existing text style (s=14) by copying
# it from an already-text cell in the same row.
# ---------------------------------------------------------------------
$ws.Range("D22").Copy($ws.Range("C23"))
$ws.Range("L23").Value = 22.222222222222

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 44
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = 15.789473684210
$ws.Range("F24").Value = 176
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = 5.389221556886
$ws.Range("I24").Value = 989
$ws.Range("J24").Value = 924
$ws.Range("K24").Value = 7.034632034632
$ws.Range("L24").Value = 26.957637997432
$ws.Range("M24").Value = 53.810264385692

# ---------------------------------------------------------------------
# Row 25 (Retail Theft)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 42
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = 23.529411764705
$ws.Range("F25").Value = 157
$ws.Range("G25").Value = 137
$ws.Range("H25").Value = 14.598540145985
$ws.Range("I25").Value = 893
$ws.Range("J25").Value = 831
$ws.Range("K25").Value = 7.460890493381
$ws.Range("L25").Value = 39.096573208722

# ---------------------------------------------------------------------
# Row 26 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 40.740740740740
$ws.Range("I26").Value = 222
$ws.Range("J26").Value = 169
$ws.Range("K26").Value = 31.360946745562
$ws.Range("L26").Value = 29.069767441860
$ws.Range("M26").Value = 27.586206896551

# ---------------------------------------------------------------------
# Row 27 (UCR Rape*)
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 5
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 166.666666666667
$ws.Range("L27").Value = 0

# ---------------------------------------------------------------------
# Row 28 (Other Sex Crimes)
# ---------------------------------------------------------------------
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 23
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -25.806451612903

# ---------------------------------------------------------------------
# Row 31 (Hate Crimes)
# ---------------------------------------------------------------------
$ws.Range("L31").Value = -28.571428571428
